$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")

# Insert a new row at position 8 (shifts existing rows 8+ down)
$ws.Rows.Item(8).Insert()

# Populate the new row with the new parameter
$ws.Cells.Item(8, 1).Value = "global"
$ws.Cells.Item(8, 2).Value = "deadHeadTripBeelineDistanceFactor"
$ws.Cells.Item(8, 3).Value = 5
$ws.Cells.Item(8, 5).Value = "The factor that is applied to travel the beeline distance at speed limit, if no path in the network is found for the dead head trip between two locations."

# Re-apply the AutoFilter over the now-larger data range (A1:E24)
$ws.AutoFilterMode = $false
$ws.Range("A1:E24").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=scenario_info!`$A`$1:`$E`$24"
    }
}

# Update the active selection to match the post-edit state
$ws.Range("C11").Select() | Out-Null
